$wb = $excel.ActiveWorkbook

# --- Remove the "6-layer" sheet entirely (its data is now obsolete) ---
$ws6 = $wb.Worksheets("6-layer")
$ws6.Delete() | Out-Null

# --- Work on the remaining "8-layer" sheet ---
$ws = $wb.Worksheets("8-layer")

# Make it the active sheet / tab, with the new selection at F22
$ws.Activate() | Out-Null
$ws.Range("F22").Select() | Out-Null

# Widen column D to fit the new "Material" values (column C stays as-is)
$ws.Columns.Item(4).ColumnWidth = 31.140625

# Add the "Grace GA-170LL" material to column D for each prepreg/laminate row,
# and bump the dielectric constant (Er) from 4.3 to 4.7 for those rows.
$materialRows = 3,5,7,9,11,13,15
foreach ($r in $materialRows) {
    $ws.Range("D$r").Value = "Grace GA-170LL"
    $ws.Range("E$r").Value = 4.7
}

# Update the trace-width (F) values for the affected rows
$ws.Range("F3").Value = 6
$ws.Range("F7").Value = 6
$ws.Range("F9").Value = 14
$ws.Range("F11").Value = 6
$ws.Range("F15").Value = 6

"done"
